# Add ancestry distribution data to the export files
#
# The "Scores" worksheet gains three new columns between the existing
# "Score and results match the original publication" column and the
# "FTP link" column:
#   - Ancestry Distribution (%) - Source of Variant Associations (GWAS)
#   - Ancestry Distribution (%) - Score Development/Training
#   - Ancestry Distribution (%) - PGS Evaluation
#
# Also the "Mapped Trait(s)" cells switch their list delimiter from ", " to "|".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

# Insert three new blank columns right before the old "FTP link" column (O).
$ws.Columns("O:Q").Insert()

# Match the formatting/style of the other header cells on row 1 (bold, bordered).
$ws.Range("N1").Copy()
$ws.Range("O1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New header row values (row 1) ---
$ws.Range("O1").Value = "Ancestry Distribution (%) - Source of Variant Associations (GWAS)"
$ws.Range("P1").Value = "Ancestry Distribution (%) - Score Development/Training"
$ws.Range("Q1").Value = "Ancestry Distribution (%) - PGS Evaluation"

# --- New data row values (row 2) ---
$ws.Range("O2").Value = "European:75.3|South Asian:13.6|East Asian:6|Hispanic or Latin American:2.2|African:1.7|Greater Middle Eastern:1.2"
$ws.Range("P2").Value = "European:100"
$ws.Range("Q2").Value = "European:42.5|East Asian:20|African:12.5|Hispanic or Latin American:12.5|South Asian:10|Additional Asian Ancestries:2.5"

# --- Update the delimiter used in the mapped-trait list cells ---
$ws.Range("D2").Value = "stroke|Ischemic stroke"
$ws.Range("E2").Value = "EFO_0000712|HP_0002140"

# --- Re-point the FTP-link hyperlink, which moved from O2 to R2 ---
$ftpUrl = "http://ftp.ebi.ac.uk/pub/databases/spot/pgs/scores/PGS000039/ScoringFiles/PGS000039.txt.gz"
$ws.Range("R2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("R2"), $ftpUrl) | Out-Null
$ws.Range("R2").Style = "Hyperlink"
$ws.Range("R2").Value = $ftpUrl
